$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "51.569.34"
$ws.Range("E2").Value = "  +1.45%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.988.52"

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "381.26"
$ws.Range("E5").Value = "  +3.93%  "

# Row 6 - Solana
$ws.Range("D6").Value = "104.98"
$ws.Range("E6").Value = "  +3.78%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  +1.95%  "

# Row 8 - USDC
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +3.57%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "37.47"
$ws.Range("E10").Value = "  +3.60%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.35%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +2.61%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.464.08"
$ws.Range("E13").Value = "  +3.33%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "18.49"
$ws.Range("E14").Value = "  +1.85%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "7.61"
$ws.Range("E15").Value = "  +4.09%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.996.00"
$ws.Range("E16").Value = "  +3.49%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.980"
$ws.Range("E17").Value = "  +7.11%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "51.583.67"
$ws.Range("E18").Value = "  +1.54%  "

# Row 19 - ImmutableX
$ws.Range("E19").Value = "  +4.56%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "7.47"
$ws.Range("E20").Value = "  +4.92%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").Value = "13.00"
$ws.Range("E21").Value = "  +1.90%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  +3.26%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "69.46"
$ws.Range("E23").Value = "  +2.43%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "263.19"
$ws.Range("E24").Value = "  +2.46%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "2.92"
$ws.Range("E25").Value = "  +9.97%  "

# Row 26 - Filecoin
$ws.Range("D26").Value = "8.39"
$ws.Range("E26").Value = "  +19.93%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "7.79"
$ws.Range("E27").Value = "  +27.31%  "

# Rows 28/29 - swap Hedera and Kaspa
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "0.116"
$ws.Range("E29").Value = "  +15.87%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "26.07"
$ws.Range("E30").Value = "  +2.78%  "

# Row 31 - Dai
$ws.Range("E31").Value = "  +0.02%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  +1.28%  "

# Row 33 - InjectiveProtocol
$ws.Range("D33").Value = "35.20"
$ws.Range("E33").Value = "  +4.02%  "

# Rows 34/35 - swap OKB and Toncoin
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "2.08"
$ws.Range("E34").Value = "  -1.50%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "51.02"
$ws.Range("E35").Value = "  +0.56%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.0457"
$ws.Range("E36").Value = "  +9.51%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.08%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  +3.43%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "17.25"
$ws.Range("E39").Value = "  +2.39%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +0.30%  "

# Row 41 - ARBITRUM
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  +1.83%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +4.72%  "

# Row 43 - Monero
$ws.Range("D43").Value = "125.87"
$ws.Range("E43").Value = "  +5.85%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "21.89"
$ws.Range("E44").Value = "  +1.09%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  +22.02%  "

# Row 46 - WEMIXToken
$ws.Range("D46").Value = "2.05"
$ws.Range("E46").Value = "  -1.59%  "

# Row 47 - ApeXProtocol
$ws.Range("D47").Value = "2.37"
$ws.Range("E47").Value = "  +2.39%  "

# Row 48 - Maker
$ws.Range("D48").Value = "2.043.31"
$ws.Range("E48").Value = "  +1.83%  "

# Row 49 - NEARProtocol
$ws.Range("E49").Value = "  +5.46%  "

# Row 50 - BEAM
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  +9.08%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  +4.41%  "
